$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header
$ws.Range("A1").Value = "Cluster Name"
$ws.Range("B1").Value = "Active cases"

# Update data rows 2-10 (rows 11-13 will be cleared/removed)
$data = @(
    @("3642 Fronditha Care Aged Care Clayton South", 10),
    @("3652 Regis Aged Care Dandenong North", 10),
    @("3662 Regis Brighton site", 11),
    @("4535 Fronditha Care Thornbury", 25),
    @("Confirmed Omicron Sircuit Bar Fitzroy", 10),
    @("Confirmed Omicron Variant The Peel Hotel Collingwood", 12),
    @("Diamond Valley Pork and Baxters Pork Laverton North", 10),
    @("Midfield Meat International Warrnambool", 23),
    @("Werribee Mercy Hospital Emergency Department", 12)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row++
}

# Remove old rows 11-13 that are no longer part of the table
$ws.Range("A11:B13").ClearContents()
